$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Case" (Docket ID) and "Court Order" (Court ORI / Court Order
# Issuing Date) sections -- rows 4 through 7 -- which are being replaced by
# the single new "Firearm Purchase Prohibition" abstract element.
$ws.Rows("4:7").Delete()

# Row 3 becomes the new "Firearm Purchase Prohibition" category header
# (previously held "Case").
$ws.Range("A3").Value = "Firearm Purchase Prohibition"

# Row 4 (previously row 8) keeps "Extension (code)" / "Firearm Purchase
# Prohibition Code" but the NIEM mapping now points at the new
# FirearmPurchaseProhibition/ActivityPersonAssociation structure.
$ws.Range("D4").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/fppq-res-ext:FirearmPurchaseProhibition[@structures:id=../nc:ActivityPersonAssociation/nc:Activity/@structures:ref]/me-fpp-codes:FirearmPurchaseProhibitionCode"
$ws.Rows("4").RowHeight = 45

# Subject/Person rows: the mappings move from
# j:ActivityCourtOrderAssociation/j:Subject to nc:ActivityPersonAssociation/nc:Person.
$ws.Range("D6").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonBirthDate/nc:Date"
$ws.Range("D7").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonGivenName"
$ws.Range("D8").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonMiddleName"
$ws.Range("D9").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonSurName"
$ws.Range("D10").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonNameSuffixText"
$ws.Range("D11").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/j:PersonSexCode"

# Widen column A to fit the new, longer category header text.
$ws.Columns("A").ColumnWidth = 26.1640625

# Match the author's last selection before saving.
$ws.Range("C7").Select()
